# Append new daily rows (2021-12-09 .. 2022-01-05) to the COVID tracking
# sheet, per "aggiornamento fino a 6 gennaio 2022".
#
# Columns: A = date serial, B = nuovi pos., C = somma mobile 7gg.,
#          D = somma mobile 7gg. per 100mila abitanti.
# Data previously ran through row 464 (2021-12-08); this extends it
# through row 491 (2022-01-05), matching the source diff exactly
# (note the date sequence skips serial 44549 / 2021-12-19, so Excel
# row 475 follows row 474 but jumps from 44548 to 44550 - that gap is
# intentional and present in the original data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastExistingRow = 464
$firstNewRow = 465
$lastNewRow = 491

$data = @(
    @(465, 44539, 9, 46, 139.4657854046024),
    @(466, 44540, 3, 41, 124.3064609041021),
    @(467, 44541, 2, 39, 118.242731103902),
    @(468, 44542, 7, 42, 127.3383258042022),
    @(469, 44543, 10, 34, 103.0834066034017),
    @(470, 44544, 1, 35, 106.1152715035018),
    @(471, 44545, 0, 32, 97.01967680320165),
    @(472, 44546, 7, 30, 90.95594700300154),
    @(473, 44547, 11, 38, 115.210866203802),
    @(474, 44548, 12, 48, 145.5295152048025),
    @(475, 44550, 12, 53, 160.6888397053027),
    @(476, 44551, 18, 61, 184.9437589061031),
    @(477, 44552, 2, 62, 187.9756238062032),
    @(478, 44553, 4, 66, 200.1030834066034),
    @(479, 44554, 18, 77, 233.4535973077039),
    @(480, 44555, 18, 84, 254.6766516084043),
    @(481, 44556, 3, 75, 227.3898675075039),
    @(482, 44557, 45, 108, 327.4414092108056),
    @(483, 44558, 23, 113, 342.6007337113058),
    @(484, 44559, 16, 127, 385.0468423127066),
    @(485, 44560, 14, 137, 415.3654913137071),
    @(486, 44561, 48, 167, 506.3214383167086),
    @(487, 44562, 56, 205, 621.5323045205105),
    @(488, 44563, 8, 210, 636.6916290210108),
    @(489, 44564, 69, 234, 709.4563866234121),
    @(490, 44565, 27, 238, 721.5838462238122),
    @(491, 44566, 62, 284, 861.0496316284147)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $ws.Cells.Item($r, 1).Value2 = $entry[1]
    $ws.Cells.Item($r, 2).Value2 = $entry[2]
    $ws.Cells.Item($r, 3).Value2 = $entry[3]
    $ws.Cells.Item($r, 4).Value2 = $entry[4]
}

# Column A carries the date style (style index 2: bold, bordered,
# centered, numFmt "YYYY-MM-DD HH:MM:SS") used throughout the rest of
# the column; replicate it onto the new date cells.
$ws.Range("A$lastExistingRow").Copy()
$ws.Range("A$firstNewRow`:A$lastNewRow").PasteSpecial(-4122)

$excel.CutCopyMode = 0
